$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '95.316.63'
$ws.Range("E2").Value = '  -0.92%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.580.99'
$ws.Range("E3").Value = '  -0.04%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.34'
$ws.Range("E5").Value = '  -1.49%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '653.71'
$ws.Range("E6").Value = '  +2.73%  '
$ws.Range("E7").Value = '  -0.82%  '
$ws.Range("E8").Value = '  -0.38%  '
$ws.Range("E9").Value = '  +0.10%  '
$ws.Range("E10").Value = '  -2.03%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.580.35'
$ws.Range("E11").Value = '  +0.06%  '
$ws.Range("E12").Value = '  +1.39%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '42.62'
$ws.Range("E13").Value = '  -1.17%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.50'
$ws.Range("E14").Value = '  +1.76%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.246.19'
$ws.Range("E15").Value = '  -0.36%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '95.240.17'
$ws.Range("E16").Value = '  -0.85%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000254'
$ws.Range("E17").Value = '  +0.27%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.572.06'
$ws.Range("E18").Value = '  -0.34%  '
$ws.Range("E19").Value = '  -3.38%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.56'
$ws.Range("E20").Value = '  -4.90%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.87'
$ws.Range("E21").Value = '  -1.29%  '
$ws.Range("E22").Value = '  +0.53%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '508.79'
$ws.Range("E23").Value = '  -1.17%  '
$ws.Range("E24").Value = '  -5.13%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.85'
$ws.Range("E25").Value = '  +3.11%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000196'
$ws.Range("E26").Value = '  -1.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '95.42'
$ws.Range("E27").Value = '  -1.29%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.72'
$ws.Range("E28").Value = '  +2.44%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.772.89'
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.03'
$ws.Range("E30").Value = '  -1.22%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.51'
$ws.Range("E31").Value = '  -0.09%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.143'
$ws.Range("E32").Value = '  -0.84%  '
$ws.Range("E33").Value = '  +0.43%  '
$ws.Range("E34").Value = '  +0.62%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.178'
$ws.Range("E35").Value = '  -2.44%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '31.81'
$ws.Range("E36").Value = '  +4.66%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.69'
$ws.Range("E37").Value = '  +13.37%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.65'
$ws.Range("E38").Value = '  +9.95%  '
$ws.Range("B39").Value = 'PolygonEcosystemToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.560'
$ws.Range("E39").Value = '  -1.27%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '587.47'
$ws.Range("E40").Value = '  +1.78%  '
$ws.Range("E41").Value = '  +0.04%  '
$ws.Range("E42").Value = '  -0.65%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.907'
$ws.Range("E43").Value = '  -1.83%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.82'
$ws.Range("E44").Value = '  +3.83%  '
$ws.Range("E45").Value = '  +5.78%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.73'
$ws.Range("E46").Value = '  +1.51%  '
$ws.Range("B47").Value = 'WhiteBITCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '23.40'
$ws.Range("E47").Value = '  -1.69%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '34.05'
$ws.Range("E48").Value = '  +29.74%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0415'
$ws.Range("E49").Value = '  -3.97%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.55'
$ws.Range("E50").Value = '  -0.15%  '
$ws.Range("B51").Value = 'Cosmos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.18'
$ws.Range("E51").Value = '  +0.24%  '
